{"js": "// Office.js (Word JavaScript API) script\n// Applies the changes described by the commit:\n//  1. Inserts five new paragraphs (data provisioning / feature engineering\n//     discussion) right before the \"3. Skill Normalization and Job Role\n//     Analysis\" heading.\n//  2. Extends the \"Application in This Project\" paragraph with additional\n//     sentences (and collapses its two runs into one).\n//  3. Fixes two reference URLs (Senger et al. and Boselli et al.) in the\n//     References paragraph (collapsing its runs into one).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// ------------------------------------------------------------------\n// 1) Insert the five new paragraphs before the \"3. Skill Normalization...\"\n//    heading (i.e. right after the paragraph that ends with \"...not always\n//    suitable for small applied projects.\").\n// ------------------------------------------------------------------\nlet prevPara = null; // paragraph ending \"...not always suitable for small applied projects.\"\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t.indexOf(\"3. Skill Normalization and Job Role Analysis\") === 0) {\n    prevPara = paragraphs.items[i - 1];\n    break;\n  }\n}\n\nif (!prevPara) {\n  throw new Error(\"Could not locate the paragraph preceding the '3. Skill Normalization...' heading.\");\n}\n\nconst newParagraphsText = [\n  \"Data provisioning and feature engineering are important in job-vacancy text mining, because the raw postings are noisy (HTML fragments, duplicated postings, missing fields, inconsistent titles) and models depend strongly on how text is cleaned and represented.\",\n  \"In the survey by Senger et al. (2024), many datasets and approaches rely on a skill base (e.g., ESCO or O*NET) to define labels and to standardize extracted skill mentions. This means the \\u201cdata provisioning\\u201d step is often: collect job postings, decide the granularity (whole posting, sentence, or span), and map extracted spans to a predefined taxonomy so that different spellings or synonyms end up as the same skill.\",\n  \"Tzimas et al. (2024) describe a full processing pipeline before any NLP model is applied. Their methodology includes selecting multiple reputable sources, extracting postings, then performing cleansing, normalization, and deduplication. They also highlight handling missing values as part of preprocessing, and treat information extraction (skills, occupation, employer, location, experience) as a separate final phase after data cleaning.\",\n  \"Boselli et al. (2018) and related WoLMIS work focus on collecting large-scale web vacancies from heterogeneous sources and then converting them into a consistent text classification dataset. A key feature engineering choice is to represent postings using bag-of-words / n-gram text features (often using title words and short text fields, because they are highly predictive). This representation is then used for supervised classification into a standard occupation taxonomy.\",\n  \"In the practical pipeline by SDK (2025), feature engineering is kept lightweight and interpretable: the job-description text is cleaned by removing clutter such as links and punctuation, and extracted skills are normalized so that abbreviations (e.g., \\u201cML\\u201d) map to a single canonical name. This improves clustering quality because similar roles share the same standardized skill tokens.\"\n];\n\n// Insert the five new paragraphs one after another (each insertion anchored\n// on the previous new paragraph) so they land in document order, ending\n// right before the \"3. Skill Normalization...\" heading.\nlet anchor = prevPara;\nfor (const paraText of newParagraphsText) {\n  anchor = anchor.insertParagraph(paraText, Word.InsertLocation.after);\n}\nawait context.sync();\n\n// ------------------------------------------------------------------\n// 2) Extend the \"Application in This Project\" paragraph.\n// ------------------------------------------------------------------\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet appPara = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t.indexOf(\"This project follows a practical applied approach\") === 0) {\n    appPara = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!appPara) {\n  throw new Error(\"Could not locate the 'This project follows a practical applied approach...' paragraph.\");\n}\n\nconst extendedText =\n  \"This project follows a practical applied approach. A predefined skill dictionary is used to extract skills from job descriptions collected from multiple datasets. Extracted skills are normalized to ensure consistency. In this project, a similar lightweight preprocessing and normalization approach is used. The focus is not on complex NLP pipelines, but on making the data clean, structured, and understandable, so that the models remain transparent and suitable for an applied student project where explainability and practical usability are important.\";\n\nappPara.getRange().insertText(extendedText, Word.InsertLocation.replace);\nawait context.sync();\n\n// ------------------------------------------------------------------\n// 3) Fix the two reference URLs in the References paragraph.\n// ------------------------------------------------------------------\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet refPara = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t.indexOf(\"[1] Senger\") === 0) {\n    refPara = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!refPara) {\n  throw new Error(\"Could not locate the References paragraph.\");\n}\n\nrefPara.load(\"text\");\nawait context.sync();\n\nlet refText = refPara.text;\nrefText = refText\n  .replace(\"https://aclanthology.org/2024.nlp4hr-1.3/\", \"https://aclanthology.org/2024.nlp4hr-1.1/\")\n  .replace(\n    \"https://link.springer.com/article/10.1007/s10844-018-0517-6\",\n    \"https://link.springer.com/article/10.1007/s10844-017-0488-x\"\n  );\n\nrefPara.getRange().insertText(refText, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script\n# Applies the changes described by the commit:\n#  1. Inserts five new paragraphs (data provisioning / feature engineering\n#     discussion) right before the \"3. Skill Normalization and Job Role\n#     Analysis\" heading.\n#  2. Extends the \"Application in This Project\" paragraph with additional\n#     sentences (and collapses its two runs into one).\n#  3. Fixes two reference URLs (Senger et al. and Boselli et al.) in the\n#     References paragraph (collapsing its runs into one).\n\n$d = $word.ActiveDocument\n\n# ------------------------------------------------------------------\n# 1) Insert the five new paragraphs right before the \"3. Skill\n#    Normalization and Job Role Analysis\" heading.\n# ------------------------------------------------------------------\n$findRange = $d.Content\n$findRange.Find.Execute(\"3. Skill Normalization and Job Role Analysis\") | Out-Null\n$headingPara = $findRange.Paragraphs(1)\n$anchorPara = $headingPara.Previous()\n\n$newParagraphsText = @(\n  \"Data provisioning and feature engineering are important in job-vacancy text mining, because the raw postings are noisy (HTML fragments, duplicated postings, missing fields, inconsistent titles) and models depend strongly on how text is cleaned and represented.\",\n  \"In the survey by Senger et al. (2024), many datasets and approaches rely on a skill base (e.g., ESCO or O*NET) to define labels and to standardize extracted skill mentions. This means the \u201cdata provisioning\u201d step is often: collect job postings, decide the granularity (whole posting, sentence, or span), and map extracted spans to a predefined taxonomy so that different spellings or synonyms end up as the same skill.\",\n  \"Tzimas et al. (2024) describe a full processing pipeline before any NLP model is applied. Their methodology includes selecting multiple reputable sources, extracting postings, then performing cleansing, normalization, and deduplication. They also highlight handling missing values as part of preprocessing, and treat information extraction (skills, occupation, employer, location, experience) as a separate final phase after data cleaning.\",\n  \"Boselli et al. (2018) and related WoLMIS work focus on collecting large-scale web vacancies from heterogeneous sources and then converting them into a consistent text classification dataset. A key feature engineering choice is to represent postings using bag-of-words / n-gram text features (often using title words and short text fields, because they are highly predictive). This representation is then used for supervised classification into a standard occupation taxonomy.\",\n  \"In the practical pipeline by SDK (2025), feature engineering is kept lightweight and interpretable: the job-description text is cleaned by removing clutter such as links and punctuation, and extracted skills are normalized so that abbreviations (e.g., \u201cML\u201d) map to a single canonical name. This improves clustering quality because similar roles share the same standardized skill tokens.\"\n)\n\n$rng = $anchorPara.Range\nforeach ($paraText in $newParagraphsText) {\n    $rng.InsertParagraphAfter()\n    $rng = $rng.Next(4)          # wdParagraph = 4 -> move range onto the new (empty) paragraph\n    $rng.Text = $paraText\n}\n\n# ------------------------------------------------------------------\n# 2) Extend the \"Application in This Project\" paragraph.\n# ------------------------------------------------------------------\n$appRange = $d.Content\n$appRange.Find.Execute(\"This project follows a practical applied approach\") | Out-Null\n$appPara = $appRange.Paragraphs(1)\n$appParaRange = $appPara.Range\n$appParaRange.MoveEnd(1, -1) | Out-Null   # exclude the paragraph mark\n$appParaRange.Text = \"This project follows a practical applied approach. A predefined skill dictionary is used to extract skills from job descriptions collected from multiple datasets. Extracted skills are normalized to ensure consistency. In this project, a similar lightweight preprocessing and normalization approach is used. The focus is not on complex NLP pipelines, but on making the data clean, structured, and understandable, so that the models remain transparent and suitable for an applied student project where explainability and practical usability are important.\"\n\n# ------------------------------------------------------------------\n# 3) Fix the two reference URLs in the References paragraph.\n# ------------------------------------------------------------------\n$refFindRange = $d.Content\n$refFindRange.Find.Execute(\"[1] Senger\") | Out-Null\n$refPara = $refFindRange.Paragraphs(1)\n$refParaRange = $refPara.Range\n$refParaRange.MoveEnd(1, -1) | Out-Null   # exclude the paragraph mark\n\n$refText = $refParaRange.Text\n$refText = $refText.Replace(\"https://aclanthology.org/2024.nlp4hr-1.3/\", \"https://aclanthology.org/2024.nlp4hr-1.1/\")\n$refText = $refText.Replace(\"https://link.springer.com/article/10.1007/s10844-018-0517-6\", \"https://link.springer.com/article/10.1007/s10844-017-0488-x\")\n$refParaRange.Text = $refText\n\nWrite-Output \"done\"\n"}
